# Update "想去人数" (want-to-go count) values in column F
# on sheets "展览" and "全部类型" to match the new scraped data snapshot.

$wb = $excel.ActiveWorkbook

# Map: sheet name -> hashtable of cell address -> new value
$updates = @{
    "展览" = @{
        "F5"  = 7
        "F6"  = 195
        "F7"  = 4558
        "F9"  = 119
        "F10" = 102
        "F12" = 86
        "F13" = 689
        "F14" = 180
        "F15" = 972
        "F16" = 78
        "F20" = 112
        "F22" = 3481
        "F23" = 5823
        "F25" = 29
        "F29" = 3349
        "F32" = 2462
        "F35" = 123
        "F36" = 210
        "F38" = 348
        "F39" = 121
        "F41" = 902
        "F42" = 17
        "F43" = 19
        "F45" = 45
        "F46" = 467
        "F48" = 548
    }
    "全部类型" = @{
        "F5"  = 7
        "F6"  = 195
        "F7"  = 4558
        "F9"  = 119
        "F10" = 102
        "F13" = 86
        "F14" = 689
        "F15" = 180
        "F16" = 972
        "F17" = 78
        "F21" = 112
        "F23" = 3481
        "F24" = 5823
        "F26" = 29
        "F30" = 3349
        "F33" = 2462
        "F36" = 123
        "F37" = 210
        "F39" = 348
        "F40" = 121
        "F42" = 902
        "F43" = 17
        "F44" = 19
        "F46" = 45
        "F47" = 467
        "F49" = 548
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellMap = $updates[$sheetName]
    foreach ($addr in $cellMap.Keys) {
        $ws.Range($addr).Value = $cellMap[$addr]
    }
}
